# lineup.xlsx — "added correcter if lineup chages"
#
# Fills in missing country / genre / comment columns for a handful of bands
# (rows 146-148, 151-155), adds a new comment to the Hypocrisy row (109),
# widens/adds columns C/D/E to fit the new comment text, and moves the
# active selection down onto column E where the new comments live.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New comment text for Hypocrisy (row 109) -----------------------
$ws.Range("E109").Value = "Peter Tagtgren is my spirit animal <3"

# --- 2. Fill in genre / comment for "The Ruins of Beverast" (row 146) --
# This band's row gets promoted to the bold "headliner" style used by the
# other first-row-of-the-day entries, matching rows like 113/114.
$ws.Rows.Item(146).Font.Bold = $true
$ws.Range("D146").Value = "atmos black/doom"
$ws.Range("E146").Value = "počasn ampak dobr"

# --- 3. Fill in comment for "Impaled nazarene" (row 147) ---------------
$ws.Rows.Item(147).Font.Bold = $true
$ws.Range("E147").Value = "definitivno (korpiklaani bojo itk zamujal k bojo pjani)"

# --- 4. Fill in genre / comment for "Tiamat" (row 148) ------------------
$ws.Range("D148").Value = "gothic/doom neki"
$ws.Range("E148").Value = "preveč rocki (vsaj poznejš) za zaklučt bo dobr"

# --- 5. Fill in country / genre / comment for the "third" stage bands --
$ws.Range("C151").Value = "rus"
$ws.Range("D151").Value = "brutal death"
$ws.Range("E151").Value = "še kr uredu loh blo bolš"

$ws.Range("C152").Value = "de"
$ws.Range("D152").Value = "brutal death"
$ws.Range("E152").Value = "propr brutal"

$ws.Range("C153").Value = "au"
$ws.Range("D153").Value = "death/black"
$ws.Range("E153").Value = "ni nek presežk"

$ws.Range("C154").Value = "de"
$ws.Range("D154").Value = "melo death/doom"
$ws.Range("E154").Value = "ni slabo bi blo pa lahko bolš"

$ws.Range("C155").Value = "lux"
$ws.Range("D155").Value = "death"
$ws.Range("E155").Value = "kr klasičn death"

# --- 6. Column widths: shrink C, widen D, add a wide comment column E --
$ws.Columns.Item(3).ColumnWidth = 4.27
$ws.Columns.Item(4).ColumnWidth = 17.8
$ws.Columns.Item(5).ColumnWidth = 46.5

# --- 7. Move the view down a row and select column E (new comments) ----
$ws.Application.ActiveWindow.ScrollRow = 115
$ws.Range("E:E").Select() | Out-Null
